# Update "Forecast Comparison" sheet with a new Week_Start_Date column and
# corrected (non zero-padded) week labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column before column B (ASIN). Everything from B..I
# shifts right to C..J automatically (dimension grows to A1:J17).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Force column B (rows 2-17) to be stored as text so the ISO date strings
# are not reinterpreted as date serial numbers.
$ws.Range("B2:B17").NumberFormat = "@"

$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

$weekLabels = @(
    "W1", "W2", "W3", "W4", "W5", "W6", "W7", "W8",
    "W9", "W10", "W11", "W12", "W13", "W14", "W15", "W16"
)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weekLabels[$i]
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
}
